# Re-colour the presentation's theme (Design > Colors) from the
# "Integral" palette to the built-in "Office" palette.
#
# PowerPoint's object model exposes the twelve theme colour slots
# (Background/Text 1-2, Accent 1-6, Hyperlink, Followed Hyperlink) as
# ThemeColorScheme.Colors(1..12).RGB on any slide's (or the slide
# master's) theme. Driving every slot to the stock "Office" RGB values
# reproduces the effect of applying the default Office theme colours
# to the deck, which is what changed in ppt/theme/theme1.xml.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cs = $s.ThemeColorScheme

# Index -> (slot, new "Office" RGB as a VBA-style BGR integer)
# 1  dk1       000000 -> 0
# 2  lt1       FFFFFF -> 16777215
# 3  dk2       44546A -> 6968388
# 4  lt2       E7E6E6 -> 15132391
# 5  accent1   5B9BD5 -> 13998939
# 6  accent2   ED7D31 -> 3243501
# 7  accent3   A5A5A5 -> 10855845
# 8  accent4   FFC000 -> 49407
# 9  accent5   4472C4 -> 12874308
# 10 accent6   70AD47 -> 4697456
# 11 hlink     0563C1 -> 12673797
# 12 folHlink  954F72 -> 7491477

$cs.Colors(1).RGB = 0
$cs.Colors(2).RGB = 16777215
$cs.Colors(3).RGB = 6968388
$cs.Colors(4).RGB = 15132391
$cs.Colors(5).RGB = 13998939
$cs.Colors(6).RGB = 3243501
$cs.Colors(7).RGB = 10855845
$cs.Colors(8).RGB = 49407
$cs.Colors(9).RGB = 12874308
$cs.Colors(10).RGB = 4697456
$cs.Colors(11).RGB = 12673797
$cs.Colors(12).RGB = 7491477
